$wb = $excel.ActiveWorkbook

# --- "About" sheet ---
$about = $wb.Worksheets.Item("About")

# Update the "Source:" value from a Google unit-converter link to "none needed"
$about.Range("B4").Value = "none needed"

# Remove the hyperlink that lived on B5 before deleting the now-empty row
if ($about.Range("B5").Hyperlinks.Count -gt 0) {
    $about.Range("B5").Hyperlinks.Delete()
}

# Delete row 5 entirely (the blank A5 / hyperlinked B5 row), shifting rows 6+ up
$about.Rows.Item(5).Delete()

# After the row shift, the unit labels move from rows 13/14 to rows 12/13
$about.Range("A12").Value = "trillion passenger-miles"
$about.Range("A13").Value = "trillion freight ton-miles"

# --- "CDCF-PMpPDOU" sheet ---
$pm = $wb.Worksheets.Item("CDCF-PMpPDOU")
$pm.Range("B2").Formula = "=10^12"

# --- "CDCF-FTMpFDOU" sheet ---
$ftm = $wb.Worksheets.Item("CDCF-FTMpFDOU")
$ftm.Range("B2").Formula = "=10^12"
